$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Bug List")
$ws2 = $wb.Worksheets.Item("Enhancements")
$ws3 = $wb.Worksheets.Item("Versions")

# --- Enhancements sheet: new row 6, plus "Added in Version" values on rows 4-6 ---
# Write the brand-new label text first so it lands earlier in the shared-string table.
$ws2.Range("A6").Value = "Add in the label printing to the delivery sheet"

# Date cell for the new row: copy the existing date cell's format (m/d/yyyy) then set the value.
$ws2.Range("B4").Copy()
$ws2.Range("B6").PasteSpecial(-4122)
$ws2.Range("B6").Value = 42499

# Version numbers added in (now that the new row's label string already exists).
$ws2.Range("C4").Value = "1.0.5"
$ws2.Range("C5").Value = "1.0.5"
$ws2.Range("C6").Value = "1.0.5"

# --- Versions sheet: new row 7 for the 1.0.5 release ---
$ws3.Range("A7").Value = "1.0.5"
$ws3.Range("C7").Value = "Bug fixes (See Bug List) and adding the labels"

# --- Selection / active-cell bookkeeping (matches final cursor positions) ---
$ws1.Range("A8").Select()
$ws3.Range("C7").Select()
$ws2.Range("A6").Select()
